# Add time measurements to output file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 - new measurement (establishes shared-string order: 135,136,137)
$ws.Range("P19").Value = "55313,4…"
$ws.Range("Q19").Value = "22000s"
$ws.Range("R19").Value = "half non elementary"

# Row 20 - new measurement (shares "half non elementary" string; adds 138,139)
$ws.Range("P20").Value = "28501,73…"
$ws.Range("Q20").Value = "5200s"
$ws.Range("R20").Value = "half non elementary"

# Row 17 - new measurement (adds 140,141)
$ws.Range("P17").Value = "12568.6…"
$ws.Range("Q17").Value = "4600s"

# Row 11 - timeout / time limit columns + new measurement (adds 142,143)
$ws.Range("I11").Value = "-"
$ws.Range("J11").Value = "-"
$ws.Range("K11").Value = "tl=86400"
$ws.Range("P11").Value = "787.0…"
$ws.Range("Q11").Value = "7300s"

# Row 21 - new measurement (adds 144,145)
$ws.Range("P21").Value = "45271.2…"
$ws.Range("Q21").Value = "14000s"

# Row 9 - timeout / time limit columns (reuses existing strings 113,114)
$ws.Range("I9").Value = "-"
$ws.Range("J9").Value = "-"
$ws.Range("K9").Value = "tl=86400"

# Row 14 - timeout / time limit columns (reuses existing strings 113,114)
$ws.Range("I14").Value = "-"
$ws.Range("J14").Value = "-"
